$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new weekly progress row (row 36) that was added for 13/6/2025.
$ws.Range("D36").Value = "13/6/2025"
$ws.Range("D36").HorizontalAlignment = -4152   # xlRight

$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 518
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 1012
$ws.Range("E36:I36").HorizontalAlignment = -4152   # xlRight

$ws.Range("J36").Value = "N/A"
$ws.Range("J36").HorizontalAlignment = -4131   # xlLeft

# Update the view: zoom out to 55% and move the active selection to H38.
$excel.ActiveWindow.Zoom = 55
$ws.Range("H38").Select()
